$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd step-id ordering ("1001-1002-1009-1005" -> "1001-1009-1002-1005").
# All three rows in column J share this value, so update them all so the
# underlying shared string is corrected everywhere it is used.
$oldValue = "1001-1002-1009-1005"
$newValue = "1001-1009-1002-1005"

$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    for ($c = 1; $c -le $used.Columns.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq $oldValue) {
            $cell.Value = $newValue
        }
    }
}

# Move the active selection from J5 to J3, and scroll the view back so the
# frozen/left-most visible column resets to A (matches the saved view state).
[void]$ws.Range("J3").Select()

